$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.402.54'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '1.633.93'
$ws.Range("E3").Value = '  -1.33%  '
$ws.Range("E4").Value = '  +0.58%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.33'
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3781'
$ws.Range("E7").Value = '  +0.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '51.99'
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3619'
$ws.Range("E9").Value = '  -1.09%  '
$ws.Range("E10").Value = '  +0.20%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.225'
$ws.Range("E11").Value = '  -2.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("E13").Value = '  -3.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.471'
$ws.Range("E14").Value = '  -3.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.347'
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001240'
$ws.Range("E16").Value = '  -2.77%  '
$ws.Range("D17").Value = '1.631.76'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.85'
$ws.Range("E18").Value = '  -0.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06929'
$ws.Range("E19").Value = '  +0.51%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.571'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.50'
$ws.Range("E21").Value = '  -4.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.51'
$ws.Range("E23").Value = '  -2.91%  '
$ws.Range("D24").Value = '23.412.16'
$ws.Range("E24").Value = '  -0.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.494'
$ws.Range("E25").Value = '  +3.30%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.062'
$ws.Range("E26").Value = '  -3.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.13'
$ws.Range("E27").Value = '  -1.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.03'
$ws.Range("E28").Value = '  +0.03%  '
$ws.Range("E29").Value = '  -0.77%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.04'
$ws.Range("E30").Value = '  -2.61%  '
$ws.Range("D31").Value = '1.812.38'
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.624'
$ws.Range("E32").Value = '  -3.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.159'
$ws.Range("E33").Value = '  -5.68%  '
$ws.Range("E34").Value = '  +8.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.38'
$ws.Range("E35").Value = '  +7.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02759'
$ws.Range("E36").Value = '  -2.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.08765'
$ws.Range("E37").Value = '  -1.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2488'
$ws.Range("E38").Value = '  -2.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.07098'
$ws.Range("E39").Value = '  -3.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.998'
$ws.Range("E40").Value = '  -4.47%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6995'
$ws.Range("E41").Value = '  -2.03%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.339'
$ws.Range("E42").Value = '  -2.76%  '
$ws.Range("E43").Value = '  -2.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '12.10'
$ws.Range("E44").Value = '  -3.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6481'
$ws.Range("E45").Value = '  -1.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.270'
$ws.Range("E47").Value = '  -3.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.968'
$ws.Range("E48").Value = '  -1.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07969'
$ws.Range("E49").Value = '  -0.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '126.51'
$ws.Range("E50").Value = '  -2.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.184'
$ws.Range("E51").Value = '  -2.78%  '
